$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.213.62"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "1.853.78"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4657"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3707"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07296"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8919"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07874"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "1.816.12"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.402"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.513"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("E17").Value = "  -0.48%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008918"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "27.238.34"
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.084"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "2.050.73"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.005"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.042"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "115.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.041"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08811"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.140"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7667"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.92%  "
$ws.Range("E34").Value = "  +3.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.519"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.83%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.705"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.109"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.02%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05218"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.943"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.044"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5113"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1625"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.492"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("E47").Value = "  -0.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.45"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("E49").Value = "  +1.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06194"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "65.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.71%  "
